$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns remain text so numeric-looking strings (trailing zeros, dot-thousands) survive
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.880.33"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").Value = "1.791.28"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "310.33"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "0.5133"
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("D8").Value = "0.3897"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "0.07833"
$ws.Range("E9").Value = "  -6.66%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.093"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "40.93"
$ws.Range("E11").Value = "  -2.44%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "6.230"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("B13").Value = "BinanceUSD"
$ws.Range("C13").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D13").Value = "0.9993"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "20.21"
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.223"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.778.66"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").Value = "91.64"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "0.06520"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "17.08"
$ws.Range("E21").Value = "  -3.96%  "
$ws.Range("D22").Value = "5.927"
$ws.Range("E22").Value = "  -2.55%  "
$ws.Range("D23").Value = "27.944.44"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("D25").Value = "2.228"
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("D26").Value = "160.34"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "20.29"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D28").Value = "1.991.38"
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").Value = "2.353"
$ws.Range("E29").Value = "  -2.74%  "
$ws.Range("D30").Value = "124.58"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").Value = "0.1075"
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("E32").Value = "  -5.45%  "
$ws.Range("D33").Value = "3.608"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").Value = "5.488"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("D35").Value = "0.07063"
$ws.Range("E35").Value = "  -6.88%  "
$ws.Range("D36").Value = "0.02304"
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("D37").Value = "8.719"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("D39").Value = "11.53"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("D40").Value = "5.009"
$ws.Range("E40").Value = "  -4.86%  "
$ws.Range("D41").Value = "0.6098"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "1.151"
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("D44").Value = "13.12"
$ws.Range("E44").Value = "  -3.46%  "
$ws.Range("D45").Value = "1.306"
$ws.Range("E45").Value = "  -6.83%  "
$ws.Range("D46").Value = "0.5907"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").Value = "3.698"
$ws.Range("D48").Value = "124.50"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "1.204"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "1.913"
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("D51").Value = "0.06817"
$ws.Range("E51").Value = "  -2.47%  "
